# Updating the Staging testdata
#
# Column K ("ExpectedFilenames") previously listed Excel/Word report-name
# templates for both "Pfizer - MM Maintenance" and "Takeda - MM Maintenance"
# products. The test data now targets Takeda only, using a new
# Standard/Excel/Word naming triple (with a "2023_" suffix on the Standard
# variant) for each of the four report sections: Clinical, Economic,
# Quality of Life and Real-world Evidence.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value  = "StandardExcelReport-Takeda - MM Maintenance-Clinical-2023_"
$ws.Range("K3").Value  = "ExcelReport-Takeda-MM Maintenance-Clinical-"
$ws.Range("K4").Value  = "WordReport-Takeda - MM Maintenance-Clinical-"
$ws.Range("K5").Value  = "StandardExcelReport-Takeda - MM Maintenance-Economic-2023_"
$ws.Range("K6").Value  = "ExcelReport-Takeda-MM Maintenance-Economic-"
$ws.Range("K7").Value  = "WordReport-Takeda - MM Maintenance-Economic-"
$ws.Range("K8").Value  = "StandardExcelReport-Takeda - MM Maintenance-Quality of Life-2023_"
$ws.Range("K9").Value  = "ExcelReport-Takeda-MM Maintenance-Quality of Life-"
$ws.Range("K10").Value = "WordReport-Takeda - MM Maintenance-Quality of Life-"
$ws.Range("K11").Value = "StandardExcelReport-Takeda - MM Maintenance-Real-world Evidence-2023_"
$ws.Range("K12").Value = "ExcelReport-Takeda-MM Maintenance-Real-world Evidence-"
$ws.Range("K13").Value = "WordReport-Takeda - MM Maintenance-Real-world Evidence-"

# The leftover Pfizer rows (previously K14:K18) are no longer needed.
$ws.Range("K14:K18").ClearContents()

# Reflect the new view/selection state: scrolled over to column I and a
# single active cell at K11 (instead of the old J1:J5 selection).
$ws.Application.ActiveWindow.ScrollColumn = 9
$ws.Range("K11").Select()
